$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2")
$rng.NumberFormat = "@"
$rng.Value = '61.441.68'
$rng.ClearFormats()

$rng = $ws.Range("E2")
$rng.NumberFormat = "@"
$rng.Value = '  -1.78%  '
$rng.ClearFormats()

$rng = $ws.Range("D3")
$rng.NumberFormat = "@"
$rng.Value = '2.991.16'
$rng.ClearFormats()

$rng = $ws.Range("E3")
$rng.NumberFormat = "@"
$rng.Value = '  -0.89%  '
$rng.ClearFormats()

$rng = $ws.Range("D4")
$rng.NumberFormat = "@"
$rng.Value = '0.999'
$rng.ClearFormats()

$rng = $ws.Range("E4")
$rng.NumberFormat = "@"
$rng.Value = '  -0.09%  '
$rng.ClearFormats()

$rng = $ws.Range("D5")
$rng.NumberFormat = "@"
$rng.Value = '589.07'
$rng.ClearFormats()

$rng = $ws.Range("E5")
$rng.NumberFormat = "@"
$rng.Value = '  +1.73%  '
$rng.ClearFormats()

$rng = $ws.Range("D6")
$rng.NumberFormat = "@"
$rng.Value = '143.62'
$rng.ClearFormats()

$rng = $ws.Range("E6")
$rng.NumberFormat = "@"
$rng.Value = '  -3.57%  '
$rng.ClearFormats()

$rng = $ws.Range("E7")
$rng.NumberFormat = "@"
$rng.Value = '  +0.06%  '
$rng.ClearFormats()

$rng = $ws.Range("D8")
$rng.NumberFormat = "@"
$rng.Value = '0.521'
$rng.ClearFormats()

$rng = $ws.Range("E8")
$rng.NumberFormat = "@"
$rng.Value = '  -0.45%  '
$rng.ClearFormats()

$rng = $ws.Range("D9")
$rng.NumberFormat = "@"
$rng.Value = '2.987.48'
$rng.ClearFormats()

$rng = $ws.Range("E9")
$rng.NumberFormat = "@"
$rng.Value = '  -1.07%  '
$rng.ClearFormats()

$rng = $ws.Range("D10")
$rng.NumberFormat = "@"
$rng.Value = '0.146'
$rng.ClearFormats()

$rng = $ws.Range("E10")
$rng.NumberFormat = "@"
$rng.Value = '  -3.29%  '
$rng.ClearFormats()

$rng = $ws.Range("D11")
$rng.NumberFormat = "@"
$rng.Value = '5.90'
$rng.ClearFormats()

$rng = $ws.Range("E11")
$rng.NumberFormat = "@"
$rng.Value = '  +4.20%  '
$rng.ClearFormats()

$rng = $ws.Range("D12")
$rng.NumberFormat = "@"
$rng.Value = '0.463'
$rng.ClearFormats()

$rng = $ws.Range("E12")
$rng.NumberFormat = "@"
$rng.Value = '  +4.79%  '
$rng.ClearFormats()

$rng = $ws.Range("D13")
$rng.NumberFormat = "@"
$rng.Value = '0.0000227'
$rng.ClearFormats()

$rng = $ws.Range("E13")
$rng.NumberFormat = "@"
$rng.Value = '  -1.51%  '
$rng.ClearFormats()

$rng = $ws.Range("D14")
$rng.NumberFormat = "@"
$rng.Value = '34.27'
$rng.ClearFormats()

$rng = $ws.Range("E14")
$rng.NumberFormat = "@"
$rng.Value = '  -2.85%  '
$rng.ClearFormats()

$rng = $ws.Range("E15")
$rng.NumberFormat = "@"
$rng.Value = '  +1.50%  '
$rng.ClearFormats()

$rng = $ws.Range("D16")
$rng.NumberFormat = "@"
$rng.Value = '3.483.26'
$rng.ClearFormats()

$rng = $ws.Range("E16")
$rng.NumberFormat = "@"
$rng.Value = '  -1.10%  '
$rng.ClearFormats()

$rng = $ws.Range("D17")
$rng.NumberFormat = "@"
$rng.Value = '7.05'
$rng.ClearFormats()

$rng = $ws.Range("E17")
$rng.NumberFormat = "@"
$rng.Value = '  +0.68%  '
$rng.ClearFormats()

$rng = $ws.Range("D18")
$rng.NumberFormat = "@"
$rng.Value = '61.390.01'
$rng.ClearFormats()

$rng = $ws.Range("E18")
$rng.NumberFormat = "@"
$rng.Value = '  -1.85%  '
$rng.ClearFormats()

$rng = $ws.Range("D19")
$rng.NumberFormat = "@"
$rng.Value = '2.987.49'
$rng.ClearFormats()

$rng = $ws.Range("E19")
$rng.NumberFormat = "@"
$rng.Value = '  -1.06%  '
$rng.ClearFormats()

$rng = $ws.Range("D20")
$rng.NumberFormat = "@"
$rng.Value = '452.56'
$rng.ClearFormats()

$rng = $ws.Range("E20")
$rng.NumberFormat = "@"
$rng.Value = '  -3.10%  '
$rng.ClearFormats()

$rng = $ws.Range("D21")
$rng.NumberFormat = "@"
$rng.Value = '14.06'
$rng.ClearFormats()

$rng = $ws.Range("E21")
$rng.NumberFormat = "@"
$rng.Value = '  +0.56%  '
$rng.ClearFormats()

$rng = $ws.Range("D22")
$rng.NumberFormat = "@"
$rng.Value = '0.686'
$rng.ClearFormats()

$rng = $ws.Range("E22")
$rng.NumberFormat = "@"
$rng.Value = '  -0.60%  '
$rng.ClearFormats()

$rng = $ws.Range("D23")
$rng.NumberFormat = "@"
$rng.Value = '7.36'
$rng.ClearFormats()

$rng = $ws.Range("E23")
$rng.NumberFormat = "@"
$rng.Value = '  -0.30%  '
$rng.ClearFormats()

$rng = $ws.Range("D24")
$rng.NumberFormat = "@"
$rng.Value = '81.81'
$rng.ClearFormats()

$rng = $ws.Range("E24")
$rng.NumberFormat = "@"
$rng.Value = '  +1.28%  '
$rng.ClearFormats()

$rng = $ws.Range("B25")
$rng.NumberFormat = "@"
$rng.Value = 'Fetch.AI'
$rng.ClearFormats()

$rng = $ws.Range("C25")
$rng.NumberFormat = "@"
$rng.Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$rng.ClearFormats()

$rng = $ws.Range("D25")
$rng.NumberFormat = "@"
$rng.Value = '2.17'
$rng.ClearFormats()

$rng = $ws.Range("E25")
$rng.NumberFormat = "@"
$rng.Value = '  -8.46%  '
$rng.ClearFormats()

$rng = $ws.Range("B26")
$rng.NumberFormat = "@"
$rng.Value = 'InternetComputer(DFINITY)'
$rng.ClearFormats()

$rng = $ws.Range("C26")
$rng.NumberFormat = "@"
$rng.Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$rng.ClearFormats()

$rng = $ws.Range("D26")
$rng.NumberFormat = "@"
$rng.Value = '12.05'
$rng.ClearFormats()

$rng = $ws.Range("E26")
$rng.NumberFormat = "@"
$rng.Value = '  -3.33%  '
$rng.ClearFormats()

$rng = $ws.Range("B27")
$rng.NumberFormat = "@"
$rng.Value = 'RenderToken'
$rng.ClearFormats()

$rng = $ws.Range("C27")
$rng.NumberFormat = "@"
$rng.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$rng.ClearFormats()

$rng = $ws.Range("D27")
$rng.NumberFormat = "@"
$rng.Value = '10.23'
$rng.ClearFormats()

$rng = $ws.Range("E27")
$rng.NumberFormat = "@"
$rng.Value = '  -2.62%  '
$rng.ClearFormats()

$rng = $ws.Range("B28")
$rng.NumberFormat = "@"
$rng.Value = 'Dai'
$rng.ClearFormats()

$rng = $ws.Range("C28")
$rng.NumberFormat = "@"
$rng.Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$rng.ClearFormats()

$rng = $ws.Range("D28")
$rng.NumberFormat = "@"
$rng.Value = '1.00'
$rng.ClearFormats()

$rng = $ws.Range("E28")
$rng.NumberFormat = "@"
$rng.Value = '  +0.13%  '
$rng.ClearFormats()

$rng = $ws.Range("B29")
$rng.NumberFormat = "@"
$rng.Value = 'PancakeSwap'
$rng.ClearFormats()

$rng = $ws.Range("C29")
$rng.NumberFormat = "@"
$rng.Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$rng.ClearFormats()

$rng = $ws.Range("D29")
$rng.NumberFormat = "@"
$rng.Value = '2.66'
$rng.ClearFormats()

$rng = $ws.Range("E29")
$rng.NumberFormat = "@"
$rng.Value = '  +1.54%  '
$rng.ClearFormats()

$rng = $ws.Range("B30")
$rng.NumberFormat = "@"
$rng.Value = 'FirstDigitalUSD'
$rng.ClearFormats()

$rng = $ws.Range("C30")
$rng.NumberFormat = "@"
$rng.Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$rng.ClearFormats()

$rng = $ws.Range("D30")
$rng.NumberFormat = "@"
$rng.Value = '1.00'
$rng.ClearFormats()

$rng = $ws.Range("E30")
$rng.NumberFormat = "@"
$rng.Value = '  -0.09%  '
$rng.ClearFormats()

$rng = $ws.Range("D31")
$rng.NumberFormat = "@"
$rng.Value = '6.96'
$rng.ClearFormats()

$rng = $ws.Range("E31")
$rng.NumberFormat = "@"
$rng.Value = '  -3.55%  '
$rng.ClearFormats()

$rng = $ws.Range("D32")
$rng.NumberFormat = "@"
$rng.Value = '2.05'
$rng.ClearFormats()

$rng = $ws.Range("E32")
$rng.NumberFormat = "@"
$rng.Value = '  -5.25%  '
$rng.ClearFormats()

$rng = $ws.Range("E33")
$rng.NumberFormat = "@"
$rng.Value = '  -0.50%  '
$rng.ClearFormats()

$rng = $ws.Range("D34")
$rng.NumberFormat = "@"
$rng.Value = '0.107'
$rng.ClearFormats()

$rng = $ws.Range("E34")
$rng.NumberFormat = "@"
$rng.Value = '  -0.72%  '
$rng.ClearFormats()

$rng = $ws.Range("D35")
$rng.NumberFormat = "@"
$rng.Value = '0.0₃0806'
$rng.ClearFormats()

$rng = $ws.Range("E35")
$rng.NumberFormat = "@"
$rng.Value = '  +1.30%  '
$rng.ClearFormats()

$rng = $ws.Range("E36")
$rng.NumberFormat = "@"
$rng.Value = '  -2.29%  '
$rng.ClearFormats()

$rng = $ws.Range("D37")
$rng.NumberFormat = "@"
$rng.Value = '5.74'
$rng.ClearFormats()

$rng = $ws.Range("E37")
$rng.NumberFormat = "@"
$rng.Value = '  -0.53%  '
$rng.ClearFormats()

$rng = $ws.Range("B38")
$rng.NumberFormat = "@"
$rng.Value = 'Stacks'
$rng.ClearFormats()

$rng = $ws.Range("C38")
$rng.NumberFormat = "@"
$rng.Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$rng.ClearFormats()

$rng = $ws.Range("D38")
$rng.NumberFormat = "@"
$rng.Value = '2.08'
$rng.ClearFormats()

$rng = $ws.Range("E38")
$rng.NumberFormat = "@"
$rng.Value = '  -3.29%  '
$rng.ClearFormats()

$rng = $ws.Range("B39")
$rng.NumberFormat = "@"
$rng.Value = 'Cosmos'
$rng.ClearFormats()

$rng = $ws.Range("C39")
$rng.NumberFormat = "@"
$rng.Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$rng.ClearFormats()

$rng = $ws.Range("D39")
$rng.NumberFormat = "@"
$rng.Value = '9.17'
$rng.ClearFormats()

$rng = $ws.Range("E39")
$rng.NumberFormat = "@"
$rng.Value = '  +2.02%  '
$rng.ClearFormats()

$rng = $ws.Range("B40")
$rng.NumberFormat = "@"
$rng.Value = 'OKB'
$rng.ClearFormats()

$rng = $ws.Range("C40")
$rng.NumberFormat = "@"
$rng.Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$rng.ClearFormats()

$rng = $ws.Range("D40")
$rng.NumberFormat = "@"
$rng.Value = '50.11'
$rng.ClearFormats()

$rng = $ws.Range("E40")
$rng.NumberFormat = "@"
$rng.Value = '  -0.15%  '
$rng.ClearFormats()

$rng = $ws.Range("D41")
$rng.NumberFormat = "@"
$rng.Value = '0.121'
$rng.ClearFormats()

$rng = $ws.Range("E41")
$rng.NumberFormat = "@"
$rng.Value = '  +8.93%  '
$rng.ClearFormats()

$rng = $ws.Range("D42")
$rng.NumberFormat = "@"
$rng.Value = '2.83'
$rng.ClearFormats()

$rng = $ws.Range("E42")
$rng.NumberFormat = "@"
$rng.Value = '  -4.03%  '
$rng.ClearFormats()

$rng = $ws.Range("D43")
$rng.NumberFormat = "@"
$rng.Value = '403.55'
$rng.ClearFormats()

$rng = $ws.Range("E43")
$rng.NumberFormat = "@"
$rng.Value = '  -4.61%  '
$rng.ClearFormats()

$rng = $ws.Range("E44")
$rng.NumberFormat = "@"
$rng.Value = '  -0.93%  '
$rng.ClearFormats()

$rng = $ws.Range("B45")
$rng.NumberFormat = "@"
$rng.Value = 'Maker'
$rng.ClearFormats()

$rng = $ws.Range("C45")
$rng.NumberFormat = "@"
$rng.Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$rng.ClearFormats()

$rng = $ws.Range("D45")
$rng.NumberFormat = "@"
$rng.Value = '2.712.96'
$rng.ClearFormats()

$rng = $ws.Range("E45")
$rng.NumberFormat = "@"
$rng.Value = '  -3.37%  '
$rng.ClearFormats()

$rng = $ws.Range("B46")
$rng.NumberFormat = "@"
$rng.Value = 'TheGraph'
$rng.ClearFormats()

$rng = $ws.Range("C46")
$rng.NumberFormat = "@"
$rng.Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$rng.ClearFormats()

$rng = $ws.Range("D46")
$rng.NumberFormat = "@"
$rng.Value = '0.265'
$rng.ClearFormats()

$rng = $ws.Range("E46")
$rng.NumberFormat = "@"
$rng.Value = '  -5.43%  '
$rng.ClearFormats()

$rng = $ws.Range("D47")
$rng.NumberFormat = "@"
$rng.Value = '38.26'
$rng.ClearFormats()

$rng = $ws.Range("E47")
$rng.NumberFormat = "@"
$rng.Value = '  +0.35%  '
$rng.ClearFormats()

$rng = $ws.Range("D48")
$rng.NumberFormat = "@"
$rng.Value = '132.73'
$rng.ClearFormats()

$rng = $ws.Range("E48")
$rng.NumberFormat = "@"
$rng.Value = '  +3.26%  '
$rng.ClearFormats()

$rng = $ws.Range("E49")
$rng.NumberFormat = "@"
$rng.Value = '  +0.10%  '
$rng.ClearFormats()

$rng = $ws.Range("E50")
$rng.NumberFormat = "@"
$rng.Value = '  -0.75%  '
$rng.ClearFormats()

$rng = $ws.Range("D51")
$rng.NumberFormat = "@"
$rng.Value = '2.16'
$rng.ClearFormats()

$rng = $ws.Range("E51")
$rng.NumberFormat = "@"
$rng.Value = '  +1.05%  '
$rng.ClearFormats()

